# Add a second slide ("Title and Content" layout == index 2) right after
# the existing slide, give it a title of "Test1" (typed with an en-GB
# keyboard while the document's own default language stays hu-HU), and
# leave its content placeholder empty -- mirroring the author's upload
# that introduced ppt/slides/slide2.xml and the new <p:sldId id="257".../>
# entry in presentation.xml.

$p = $ppt.ActivePresentation

# ppLayoutObject (2) = "Title and Content" slide layout, inserted as slide 2.
$s = $p.Slides.Add(2, 2)

$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Test1"
$title.LanguageID = "en-GB"

Write-Output "Added slide 2 with title 'Test1'"
